# Updated cryptos list snapshot values (price + volume/1h) pulled by the
# scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, percentage strings, and price
# strings that already contain more than one "." so Excel cannot treat
# them as a number) -- these can be written directly.
$textUpdates = @{
    "D2" = "61.868.89"
    "E2" = "  -1.77%  "
    "D3" = "2.913.26"
    "E3" = "  -2.29%  "
    "E4" = "  +0.17%  "
    "E5" = "  -1.57%  "
    "E6" = "  +0.50%  "
    "E7" = "  +0.06%  "
    "E8" = "  +0.52%  "
    "D9" = "2.913.22"
    "E9" = "  -2.20%  "
    "E10" = "  -5.03%  "
    "E11" = "  +4.11%  "
    "E12" = "  -3.14%  "
    "E13" = "  +1.40%  "
    "E14" = "  -1.70%  "
    "E15" = "  -1.68%  "
    "D16" = "3.399.17"
    "E16" = "  -2.03%  "
    "D17" = "61.915.99"
    "E17" = "  -1.45%  "
    "E18" = "  -2.16%  "
    "D19" = "2.912.14"
    "E19" = "  -1.23%  "
    "E20" = "  -2.07%  "
    "E21" = "  -0.95%  "
    "E22" = "  -2.98%  "
    "E23" = "  -3.34%  "
    "E24" = "  -1.46%  "
    "E25" = "  -2.29%  "
    "E26" = "  -6.02%  "
    "E27" = "  -3.82%  "
    "E28" = "  +0.04%  "
    "E29" = "  +21.86%  "
    "E30" = "  +1.35%  "
    "E31" = "  -2.24%  "
    "E32" = "  -1.07%  "
    "E33" = "  +0.27%  "
    "B34" = "FirstDigitalUSD"
    "C34" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "E34" = "  -0.02%  "
    "B35" = "EthereumClassic"
    "C35" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "E35" = "  -2.67%  "
    "E36" = "  -1.68%  "
    "E37" = "  +3.98%  "
    "E38" = "  -2.69%  "
    "B39" = "OKB"
    "C39" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "E39" = "  -1.04%  "
    "B40" = "Stacks"
    "C40" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "E40" = "  -2.25%  "
    "E41" = "  -3.39%  "
    "E42" = "  -2.26%  "
    "E43" = "  -5.28%  "
    "E44" = "  -0.11%  "
    "D45" = "2.700.87"
    "E45" = "  -0.26%  "
    "E46" = "  +0.21%  "
    "E47" = "  -1.63%  "
    "E48" = "  -7.52%  "
    "E49" = "  +0.04%  "
    "E50" = "  -1.52%  "
    "E51" = "  -4.64%  "
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}

# Price strings in column D that look like plain decimals (e.g. "1.00",
# "0.150") would otherwise be auto-converted to numbers by Excel and lose
# their original text formatting (trailing/leading zeros). Force the cell
# to Text format before writing, then clear the format delta again so the
# cell style stays identical to its neighbours (matches source data which
# stores every Price/Volume cell as an unstyled inline string).
$numericLookingUpdates = @{
    "D4" = "1.00"
    "D5" = "587.54"
    "D6" = "146.77"
    "D10" = "7.02"
    "D11" = "0.150"
    "D12" = "0.436"
    "D13" = "0.0000239"
    "D14" = "33.04"
    "D18" = "6.61"
    "D20" = "435.17"
    "D21" = "13.45"
    "D22" = "0.659"
    "D23" = "6.94"
    "D24" = "81.04"
    "D25" = "11.85"
    "D26" = "10.22"
    "D27" = "2.07"
    "D30" = "7.23"
    "D31" = "2.56"
    "D32" = "2.11"
    "D33" = "0.109"
    "D34" = "1.00"
    "D35" = "25.89"
    "D36" = "0.976"
    "D37" = "3.08"
    "D38" = "5.51"
    "D39" = "49.13"
    "D40" = "2.01"
    "D41" = "8.36"
    "D43" = "0.272"
    "D44" = "38.93"
    "D46" = "134.53"
    "D48" = "344.69"
    "D51" = "22.36"
}

foreach ($cellRef in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$cellRef]
    $cell.ClearFormats()
}
